$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new time-log entry as row 17 (Date, Time Spent, Descriptions)
$ws.Range("A17").Value2 = 45238
$ws.Range("B17").Value2 = "~1hr"
$ws.Range("C17").Value2 = "commenting the translator.c and .h files"

# Move the selection to where the cursor ends up after entering the row
$ws.Range("C18").Select()
